$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "_tejgtotfun_f5amb"
$ws.Range("C2").Value = 0.03002000106757184
$ws.Range("B3").Value = "_tejgfun_f2ct05pgrco"
$ws.Range("C3").Value = 0.01867467115288111
$ws.Range("B4").Value = "_tejgtotfun_f2pgrco"
$ws.Range("C4").Value = 0.01051242608725471
$ws.Range("B5").Value = "_tejgfun_f5ct05opseg"
$ws.Range("C5").Value = 0.00763376342751429
$ws.Range("B6").Value = "compu_muni_5"
$ws.Range("C6").Value = 0.006460877857838675
$ws.Range("B7").Value = "_tejgtotfun_f2opseg"
$ws.Range("C7").Value = 0.00535427581093951
$ws.Range("B8").Value = "_tejgfun_f5ct05pgrco"
$ws.Range("C8").Value = 0.005278559730775546
$ws.Range("B9").Value = "_tejgct_r09gstcr"
$ws.Range("C9").Value = 0.00428265947448246
$ws.Range("B10").Value = "dfgdevpiagfun_f5r18ct05transpc"
$ws.Range("C10").Value = 0.003542455175140605
$ws.Range("B11").Value = "_tejgtotfun_f5pgrco"
$ws.Range("C11").Value = 0.003470145944405625
$ws.Range("B12").Value = "_dfgdevpiagge_r09ct05pobso"
$ws.Range("C12").Value = 0.002537196169125774
$ws.Range("B13").Value = "_tejgrb_redr"
$ws.Range("C13").Value = 0.002489239240342244
$ws.Range("B14").Value = "_devppimfun_f2ct05ind"
$ws.Range("C14").Value = 0.002440416190272476
$ws.Range("B15").Value = "_tejgge_r09ct05biser"
$ws.Range("C15").Value = 0.002422350692369926
$ws.Range("B16").Value = "dfgpimpiafun_f5r18ct06pgrco"
$ws.Range("C16").Value = 0.002399428535304918
$ws.Range("B17").Value = "_tejgfun_f5r08ct05pgrcopc"
$ws.Range("C17").Value = 0.00208420263588348
$ws.Range("B18").Value = "dfgpimpiage_r08ct05dotra"
$ws.Range("C18").Value = 0.002071055752040122
$ws.Range("B19").Value = "_tejgft_redr"
$ws.Range("C19").Value = 0.001957159008484327
$ws.Range("B20").Value = "dfgdevpiagfun_f5r18ct06pgrco"
$ws.Range("C20").Value = 0.001953594808568985
$ws.Range("B21").Value = "dfgdevpiagge_r00ct05biser"
$ws.Range("C21").Value = 0.001896939070291455
$ws.Range("B22").Value = "dfgpimpiatotfun_f5pgrco"
$ws.Range("C22").Value = 0.001817638524801889
$ws.Range("B23").Value = "_tejgge_r08ct05biser"
$ws.Range("C23").Value = 0.001803351475839784
$ws.Range("B24").Value = "_dfgpimpiagge_r09ct05otgst"
$ws.Range("C24").Value = 0.001779235663121049
$ws.Range("B25").Value = "_dfgpimpiafun_f5ct06opsegpc"
$ws.Range("C25").Value = 0.001776720473692271
$ws.Range("B26").Value = "dfgdevpiagfun_f5r07ct05pgrcopc"
$ws.Range("C26").Value = 0.00170975753021845
$ws.Range("B27").Value = "tejgge_r07ct05otgstpc"
$ws.Range("C27").Value = 0.00170973117518319
$ws.Range("B28").Value = "dfgdevpiagfun_f5ct05sanpc"
$ws.Range("C28").Value = 0.001625716086467806
$ws.Range("B29").Value = "_tejgct_r09gstcrpc"
$ws.Range("C29").Value = 0.001622978304318528
$ws.Range("B30").Value = "dfgdevpiagfun_f5ct05transpc"
$ws.Range("C30").Value = 0.001620755641833165
$ws.Range("B31").Value = "_tejgfun_f5ct05amb"
$ws.Range("C31").Value = 0.001616406908461108
$ws.Range("B32").Value = "orgs_3"
$ws.Range("C32").Value = 0.001536330062520839
$ws.Range("B33").Value = "pimgfun_f1ct06san"
$ws.Range("C33").Value = 0.001531763891005311
$ws.Range("B34").Value = "_piagge_r09ct05otgstpc"
$ws.Range("C34").Value = 0.001511745792653283
$ws.Range("B35").Value = "_tejgge_r08ct05pobso"
$ws.Range("C35").Value = 0.001502753553425762
$ws.Range("B36").Value = "pimgfun_f5r18ct05opseg"
$ws.Range("C36").Value = 0.001482912450858794
$ws.Range("B37").Value = "pimgtotfun_f5trans"
$ws.Range("C37").Value = 0.001480073803481577
$ws.Range("B38").Value = "devppimfun_f2ct06agro"
$ws.Range("C38").Value = 0.001471117250740902
$ws.Range("B39").Value = "_dfgdevpiagct_r18gstcp"
$ws.Range("C39").Value = 0.001464363212062589
$ws.Range("B40").Value = "piagtotfun_f5r07protspc"
$ws.Range("C40").Value = 0.001442669562939713
$ws.Range("B41").Value = "dfgpimpiafun_f5ct05sanpc"
$ws.Range("C41").Value = 0.001435588685041514
$ws.Range("B42").Value = "piagge_r18ct06acanfpc"
$ws.Range("C42").Value = 0.001430916248202294
$ws.Range("B43").Value = "_dfgdevpiagge_r09ct05popso"
$ws.Range("C43").Value = 0.001417916097694246
$ws.Range("B44").Value = "dfgdevpiagtotfun_f4transpc"
$ws.Range("C44").Value = 0.001410103777947798
$ws.Range("B45").Value = "_dfgpimpiatotfun_f5edu"
$ws.Range("C45").Value = 0.001391649698358097
$ws.Range("B46").Value = "pimgfun_f5ct06trans"
$ws.Range("C46").Value = 0.001372279591767399
$ws.Range("B47").Value = "devppimtotfun_f5r18amb"
$ws.Range("C47").Value = 0.00134644840544339
$ws.Range("B48").Value = "tdvgfun_f5ct05prots"
$ws.Range("C48").Value = 0.001302789619331439
$ws.Range("B49").Value = "pimgtotfun_f1san"
$ws.Range("C49").Value = 0.001285009267155171
$ws.Range("B50").Value = "dfgdevpiagfun_f5r07ct06protspc"
$ws.Range("C50").Value = 0.001283891390710773
$ws.Range("B51").Value = "_piagfun_f1ct05trans"
$ws.Range("C51").Value = 0.001269718288626472
